$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings stored as text in the
# original workbook. Temporarily mark the whole data range as Text so
# Excel does not coerce the new values into numbers, then restore the
# default (Normal) style so the cell formatting matches the original.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.892.44'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").Value = '2.297.60'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '114.17'
$ws.Range("E5").Value = '  +18.49%  '
$ws.Range("D6").Value = '269.81'
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = '0.629'
$ws.Range("E7").Value = '  +0.58%  '
$ws.Range("E8").Value = '  +0.26%  '
$ws.Range("E9").Value = '  +1.69%  '
$ws.Range("D10").Value = '48.21'
$ws.Range("E10").Value = '  +6.54%  '
$ws.Range("D11").Value = '0.0948'
$ws.Range("E11").Value = '  +1.28%  '
$ws.Range("D12").Value = '9.05'
$ws.Range("E12").Value = '  +14.41%  '
$ws.Range("E13").Value = '  +0.05%  '
$ws.Range("D14").Value = '15.87'
$ws.Range("E14").Value = '  +1.01%  '
$ws.Range("D15").Value = '2.640.36'
$ws.Range("E15").Value = '  +0.25%  '
$ws.Range("D16").Value = '0.856'
$ws.Range("E16").Value = '  -0.32%  '
$ws.Range("D17").Value = '2.287.52'
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D18").Value = '43.785.04'
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").Value = '0.0000111'
$ws.Range("E19").Value = '  -0.63%  '
$ws.Range("D20").Value = '6.85'
$ws.Range("E20").Value = '  +10.53%  '
$ws.Range("D21").Value = '72.17'
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").Value = '2.44'
$ws.Range("E22").Value = '  -2.05%  '
$ws.Range("D23").Value = '3.01'
$ws.Range("E23").Value = '  +10.84%  '
$ws.Range("E24").Value = '  +0.25%  '
$ws.Range("E25").Value = '  +6.14%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").Value = '11.67'
$ws.Range("E27").Value = '  +2.75%  '
$ws.Range("D28").Value = '41.82'
$ws.Range("E28").Value = '  +8.20%  '
$ws.Range("D29").Value = '3.38'
$ws.Range("E29").Value = '  -2.11%  '
$ws.Range("D30").Value = '2.26'
$ws.Range("E30").Value = '  -0.79%  '
$ws.Range("D31").Value = '175.33'
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("D32").Value = '0.0934'
$ws.Range("D33").Value = '21.56'
$ws.Range("E33").Value = '  -1.32%  '
$ws.Range("D34").Value = '5.75'
$ws.Range("E34").Value = '  +5.84%  '
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("E36").Value = '  -0.36%  '
$ws.Range("D37").Value = '0.0365'
$ws.Range("E37").Value = '  +3.79%  '
$ws.Range("E38").Value = '  +0.09%  '
$ws.Range("D39").Value = '3.85'
$ws.Range("E39").Value = '  +7.18%  '
$ws.Range("D40").Value = '74.91'
$ws.Range("E40").Value = '  +16.32%  '
$ws.Range("D41").Value = '0.246'
$ws.Range("E41").Value = '  +3.89%  '
$ws.Range("D42").Value = '13.72'
$ws.Range("E42").Value = '  +11.89%  '
$ws.Range("E43").Value = '  +2.53%  '
$ws.Range("D44").Value = '6.35'
$ws.Range("E44").Value = '  +22.05%  '
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("E46").Value = '  +3.39%  '
$ws.Range("D47").Value = '8.82'
$ws.Range("E47").Value = '  +1.41%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.0996'
$ws.Range("E48").Value = '  -2.90%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = '101.81'
$ws.Range("E49").Value = '  +3.92%  '
$ws.Range("E50").Value = '  +3.27%  '
$ws.Range("D51").Value = '0.468'
$ws.Range("E51").Value = '  +5.43%  '

$ws.Range("D2:D51").Style = "Normal"
